{"js": "// Update the division-problem prompts in the worksheet table.\n// Each (old, new) pair below is a w:t run of the form \"XX\u00f7Y=\" that needs\n// to become a different \"XX\u00f7Y=\" expression, per the commit's refreshed\n// problem set.\n\nconst pairs = [\n  [\"33\u00f79=\", \"72\u00f75=\"],\n  [\"67\u00f78=\", \"12\u00f77=\"],\n  [\"10\u00f78=\", \"17\u00f79=\"],\n  [\"20\u00f75=\", \"63\u00f72=\"],\n  [\"23\u00f77=\", \"57\u00f72=\"],\n  [\"35\u00f72=\", \"62\u00f77=\"],\n  [\"77\u00f79=\", \"76\u00f76=\"],\n  [\"83\u00f75=\", \"95\u00f74=\"],\n  [\"19\u00f77=\", \"33\u00f73=\"],\n  [\"25\u00f74=\", \"67\u00f74=\"],\n  [\"12\u00f72=\", \"96\u00f76=\"],\n  [\"36\u00f72=\", \"38\u00f72=\"],\n  [\"84\u00f76=\", \"78\u00f72=\"],\n  [\"16\u00f77=\", \"45\u00f77=\"],\n  [\"37\u00f76=\", \"97\u00f73=\"],\n  [\"67\u00f75=\", \"28\u00f75=\"],\n  [\"46\u00f73=\", \"72\u00f78=\"],\n  [\"89\u00f76=\", \"13\u00f75=\"],\n  [\"51\u00f78=\", \"44\u00f76=\"],\n  [\"90\u00f75=\", \"65\u00f73=\"],\n  [\"68\u00f75=\", \"60\u00f72=\"],\n  [\"72\u00f78=\", \"60\u00f77=\"],\n  [\"58\u00f79=\", \"52\u00f79=\"],\n  [\"23\u00f76=\", \"47\u00f75=\"],\n  [\"31\u00f74=\", \"42\u00f76=\"],\n];\n\nasync function replaceExactText(oldText, newText) {\n  const body = context.document.body;\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Some \"new\" values equal other entries' \"old\" values (e.g. \"72\u00f78=\" is\n// both a source and a target). Doing the replacements directly, in\n// sequence, could let an earlier replacement's output be clobbered by a\n// later rule meant for the original text. To make every substitution\n// independent of ordering, first swap each old value for a unique\n// placeholder, then swap every placeholder for its real new value.\nfor (let i = 0; i < pairs.length; i++) {\n  const [oldText] = pairs[i];\n  await replaceExactText(oldText, `@@MIGRATE_${i}@@`);\n}\n\nfor (let i = 0; i < pairs.length; i++) {\n  const [, newText] = pairs[i];\n  await replaceExactText(`@@MIGRATE_${i}@@`, newText);\n}\n", "ps1": "# Update the division-problem prompts in the worksheet table.\n# Each (old, new) pair below is a w:t run of the form \"XX\u00f7Y=\" that needs\n# to become a different \"XX\u00f7Y=\" expression, per the commit's refreshed\n# problem set. We replace old text with new text throughout the document\n# body using Find/Replace on $d.Content.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"33\u00f79=\", \"72\u00f75=\"),\n    @(\"67\u00f78=\", \"12\u00f77=\"),\n    @(\"10\u00f78=\", \"17\u00f79=\"),\n    @(\"20\u00f75=\", \"63\u00f72=\"),\n    @(\"23\u00f77=\", \"57\u00f72=\"),\n    @(\"35\u00f72=\", \"62\u00f77=\"),\n    @(\"77\u00f79=\", \"76\u00f76=\"),\n    @(\"83\u00f75=\", \"95\u00f74=\"),\n    @(\"19\u00f77=\", \"33\u00f73=\"),\n    @(\"25\u00f74=\", \"67\u00f74=\"),\n    @(\"12\u00f72=\", \"96\u00f76=\"),\n    @(\"36\u00f72=\", \"38\u00f72=\"),\n    @(\"84\u00f76=\", \"78\u00f72=\"),\n    @(\"16\u00f77=\", \"45\u00f77=\"),\n    @(\"37\u00f76=\", \"97\u00f73=\"),\n    @(\"67\u00f75=\", \"28\u00f75=\"),\n    @(\"46\u00f73=\", \"72\u00f78=\"),\n    @(\"89\u00f76=\", \"13\u00f75=\"),\n    @(\"51\u00f78=\", \"44\u00f76=\"),\n    @(\"90\u00f75=\", \"65\u00f73=\"),\n    @(\"68\u00f75=\", \"60\u00f72=\"),\n    @(\"72\u00f78=\", \"60\u00f77=\"),\n    @(\"58\u00f79=\", \"52\u00f79=\"),\n    @(\"23\u00f76=\", \"47\u00f75=\"),\n    @(\"31\u00f74=\", \"42\u00f76=\")\n)\n\nfunction Replace-ExactText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# Some \"new\" values equal other entries' \"old\" values (e.g. \"72\u00f78=\" is\n# both a source and a target). Doing the replacements directly, in\n# sequence, could let an earlier replacement's output be clobbered by a\n# later rule meant for the original text. To make every substitution\n# independent of ordering, first swap each old value for a unique\n# placeholder, then swap every placeholder for its real new value.\n$i = 0\nforeach ($pair in $pairs) {\n    $placeholder = \"@@MIGRATE_$i@@\"\n    Replace-ExactText $pair[0] $placeholder\n    $i++\n}\n\n$i = 0\nforeach ($pair in $pairs) {\n    $placeholder = \"@@MIGRATE_$i@@\"\n    Replace-ExactText $placeholder $pair[1]\n    $i++\n}\n\nWrite-Output \"Replaced $($pairs.Count) division expressions\"\n"}
